$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted with literal dots (e.g. thousand
# separators) that must stay text, not be coerced to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.459.14"
$ws.Range("E2").Value = "  +9.86%  "
$ws.Range("D3").Value = "1.786.82"
$ws.Range("E3").Value = "  +6.65%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "337.20"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "0.9946"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "0.3791"
$ws.Range("E7").Value = "  +3.68%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3481"
$ws.Range("E8").Value = "  +7.45%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "49.43"
$ws.Range("E9").Value = "  +4.27%  "
$ws.Range("D10").Value = "1.223"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("D11").Value = "0.07691"
$ws.Range("E11").Value = "  +7.37%  "
$ws.Range("D12").Value = "0.9977"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "6.585"
$ws.Range("E13").Value = "  +8.15%  "
$ws.Range("D14").Value = "21.49"
$ws.Range("E14").Value = "  +9.23%  "
$ws.Range("D15").Value = "7.222"
$ws.Range("E15").Value = "  +8.20%  "
$ws.Range("D16").Value = "1.778.88"
$ws.Range("E16").Value = "  +6.44%  "
$ws.Range("D17").Value = "0.00001114"
$ws.Range("E17").Value = "  +6.38%  "
$ws.Range("D18").Value = "0.06766"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "84.90"
$ws.Range("E19").Value = "  +7.54%  "
$ws.Range("D20").Value = "0.9950"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  +10.57%  "
$ws.Range("D22").Value = "6.384"
$ws.Range("E22").Value = "  +7.97%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "13.16"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "27.453.11"
$ws.Range("E24").Value = "  +9.77%  "
$ws.Range("D25").Value = "2.476"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "1.500"
$ws.Range("E26").Value = "  +24.66%  "
$ws.Range("D27").Value = "2.529"
$ws.Range("E27").Value = "  +6.07%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  +7.69%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "153.14"
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").Value = "1.975.05"
$ws.Range("E30").Value = "  +6.33%  "
$ws.Range("D31").Value = "135.23"
$ws.Range("E31").Value = "  +7.23%  "
$ws.Range("D32").Value = "4.112"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "6.279"
$ws.Range("E33").Value = "  +8.22%  "
$ws.Range("D34").Value = "0.08750"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").Value = "13.42"
$ws.Range("E35").Value = "  +8.71%  "
$ws.Range("D36").Value = "1.720"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("D37").Value = "5.627"
$ws.Range("E37").Value = "  +9.07%  "
$ws.Range("D38").Value = "0.02405"
$ws.Range("E38").Value = "  +7.94%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2247"
$ws.Range("E39").Value = "  +7.49%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06453"
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("D41").Value = "0.6776"
$ws.Range("E41").Value = "  +13.55%  "
$ws.Range("D42").Value = "8.734"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("D43").Value = "1.236"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "14.66"
$ws.Range("E44").Value = "  +7.40%  "
$ws.Range("D45").Value = "0.6459"
$ws.Range("E45").Value = "  +12.66%  "
$ws.Range("D46").Value = "0.9954"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "3.971"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("D48").Value = "2.148"
$ws.Range("E48").Value = "  +9.23%  "
$ws.Range("D49").Value = "131.30"
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").Value = "0.07349"
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("D51").Value = "80.22"
$ws.Range("E51").Value = "  +7.44%  "
